# Standardize tissue names in the Samples sheet of the example workbook.
# (Mirrors commit "Standardized tissue names in examples" / "Updated tissue names".)

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Samples")

# Map of old (display) tissue names -> new standardized tissue names.
# Both D2 and D17 originally held "BAT", so both must be updated so the
# underlying shared string is fully replaced (no leftover old text).
$ws.Range("D2").Value  = "brown_adipose_tissue"
$ws.Range("D3").Value  = "brain"
$ws.Range("D4").Value  = "diaphragm"
$ws.Range("D5").Value  = "gastrocnemius"
$ws.Range("D6").Value  = "white_adipose_tissue_gonadal"
$ws.Range("D7").Value  = "heart"
$ws.Range("D8").Value  = "kidney"
$ws.Range("D9").Value  = "liver"
$ws.Range("D10").Value = "lung"
$ws.Range("D11").Value = "pancreas"
$ws.Range("D12").Value = "quadricep"
$ws.Range("D13").Value = "small_intestine"
$ws.Range("D14").Value = "soleus"
$ws.Range("D15").Value = "spleen"
$ws.Range("D16").Value = "serum_plasma_unspecified_location"
$ws.Range("D17").Value = "brown_adipose_tissue"

# Widen the Tissue column (D) now that it holds longer standardized names.
$ws.Columns.Item(4).ColumnWidth = 30

# Move the active selection to D17, matching the post-edit cursor position.
[void]$ws.Range("D17").Select()
